$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = [double]"25.65000000000057"
$ws.Range("H2").Value = [double]"1.725445188949237e-07"
$ws.Range("I2").Value = [double]"1.725445188949237e-07"
$ws.Range("L2").Value = [double]"61.2940079136619"
$ws.Range("M2").Value = "[40.17454029306674, 82.41347553425706]"
$ws.Range("N2").Value = [double]"5.293142226747705e-07"
$ws.Range("O2").Value = [double]"5.293142226747705e-07"
$ws.Range("P2").Value = [double]"1.628973968528041"
$ws.Range("Q2").Value = "[1.2012896833546556, 2.0566582537014266]"
$ws.Range("R2").Value = [double]"1.041339015017684e-09"
$ws.Range("S2").Value = [double]"1.041339015017684e-09"
$ws.Range("T2").Value = [double]"59.96031348592786"
$ws.Range("U2").Value = "[46.49171961777904, 73.42890735407667]"
$ws.Range("V2").Value = [double]"1.412958638979944e-11"
$ws.Range("W2").Value = [double]"1.412958638979944e-11"
$ws.Range("X2").Value = [double]"19.00000000000042"
$ws.Range("Y2").Value = [double]"17.25405405405444"
$ws.Range("Z2").Value = [double]"20.74594594594641"

# Row 3
$ws.Range("F3").Value = [double]"25.65000000000057"
$ws.Range("H3").Value = [double]"0.0004014341337060801"
$ws.Range("I3").Value = [double]"0.0004014341337060801"
$ws.Range("L3").Value = [double]"48.07577717778288"
$ws.Range("M3").Value = "[19.15966678983395, 76.9918875657318]"
$ws.Range("N3").Value = [double]"0.001649566000748948"
$ws.Range("O3").Value = [double]"0.001649566000748948"
$ws.Range("P3").Value = [double]"1.239026532046425"
$ws.Range("Q3").Value = "[0.5471843060306547, 1.9308687580621946]"
$ws.Range("R3").Value = [double]"0.0007729581409410269"
$ws.Range("S3").Value = [double]"0.0007729581409410269"
$ws.Range("T3").Value = [double]"65.68526379216576"
$ws.Range("U3").Value = "[49.61156307209633, 81.75896451223518]"
$ws.Range("V3").Value = [double]"1.595186205349819e-10"
$ws.Range("W3").Value = [double]"1.595186205349819e-10"
$ws.Range("X3").Value = [double]"20.59189189189235"
$ws.Range("Y3").Value = [double]"17.76756756756797"
$ws.Range("Z3").Value = [double]"23.41621621621674"

# Row 4
$ws.Range("B4").Value = [double]"1"
$ws.Range("F4").Value = [double]"25.65000000000057"
$ws.Range("H4").Value = [double]"0.0005503476294157483"
$ws.Range("I4").Value = [double]"0.0005503476294157483"
$ws.Range("L4").Value = [double]"48.56690102763248"
$ws.Range("M4").Value = "[16.56791309752581, 80.56588895773916]"
$ws.Range("N4").Value = [double]"0.003755120077692853"
$ws.Range("O4").Value = [double]"0.003755120077692853"
$ws.Range("P4").Value = [double]"0.7736053981812709"
$ws.Range("Q4").Value = "[0.1949737182408091, 1.3522370781217328]"
$ws.Range("R4").Value = [double]"0.009918379043219128"
$ws.Range("S4").Value = [double]"0.009918379043219128"
$ws.Range("T4").Value = [double]"59.24815731406017"
$ws.Range("U4").Value = "[43.13370901067552, 75.36260561744483]"
$ws.Range("V4").Value = [double]"2.563619316831023e-09"
$ws.Range("W4").Value = [double]"2.563619316831023e-09"
$ws.Range("X4").Value = [double]"22.49189189189239"
$ws.Range("Y4").Value = [double]"20.12972972973018"
$ws.Range("Z4").Value = [double]"24.8540540540546"

# Row 5
$ws.Range("B5").Value = [double]"0"
$ws.Range("F5").Value = [double]"25.65000000000057"
$ws.Range("H5").Value = [double]"6.865853916959708e-06"
$ws.Range("I5").Value = [double]"6.865853916959708e-06"
$ws.Range("L5").Value = [double]"54.01918577236902"
$ws.Range("M5").Value = "[30.064448053591647, 77.9739234911464]"
$ws.Range("N5").Value = [double]"4.154991515914475e-05"
$ws.Range("O5").Value = [double]"4.154991515914475e-05"
$ws.Range("P5").Value = [double]"0.4591316590831935"
$ws.Range("Q5").Value = "[-0.03144737390980712, 0.9497106920761942]"
$ws.Range("R5").Value = [double]"0.06589800617748121"
$ws.Range("S5").Value = [double]"0.06589800617748121"
$ws.Range("T5").Value = [double]"68.51422907107838"
$ws.Range("U5").Value = "[55.211541096600634, 81.81691704555612]"
$ws.Range("V5").Value = [double]"1.62980740014973e-13"
$ws.Range("W5").Value = [double]"1.62980740014973e-13"
$ws.Range("X5").Value = [double]"23.7756756756762"
$ws.Range("Y5").Value = [double]"21.77297297297346"
$ws.Range("Z5").Value = [double]"25.77837837837895"

# Row 6
$ws.Range("F6").Value = [double]"25.65000000000057"
$ws.Range("H6").Value = [double]"6.269233972511223e-06"
$ws.Range("I6").Value = [double]"6.269233972511223e-06"
$ws.Range("L6").Value = [double]"56.36444723531097"
$ws.Range("M6").Value = "[33.237171055260475, 79.49172341536146]"
$ws.Range("N6").Value = [double]"1.246138708665967e-05"
$ws.Range("O6").Value = [double]"1.246138708665967e-05"
$ws.Range("P6").Value = [double]"0.1195000208572692"
$ws.Range("Q6").Value = "[-0.3333421634439624, 0.5723422051585008]"
$ws.Range("R6").Value = [double]"0.5976867489866926"
$ws.Range("S6").Value = [double]"0.5976867489866926"
$ws.Range("T6").Value = [double]"70.53758820156436"
$ws.Range("U6").Value = "[56.733944381649806, 84.34123202147892]"
$ws.Range("V6").Value = [double]"2.093880624443045e-13"
$ws.Range("W6").Value = [double]"2.093880624443045e-13"
$ws.Range("X6").Value = [double]"25.16216216216272"
$ws.Range("Y6").Value = [double]"23.31351351351403"
$ws.Range("Z6").Value = [double]"27.01081081081141"

# Row 7
$ws.Range("B7").Value = [double]"0"
$ws.Range("F7").Value = [double]"25.65000000000057"
$ws.Range("H7").Value = [double]"0.0002606650729047866"
$ws.Range("I7").Value = [double]"0.0002606650729047866"
$ws.Range("L7").Value = [double]"52.84312010433631"
$ws.Range("M7").Value = "[21.6244183688446, 84.06182183982801]"
$ws.Range("N7").Value = [double]"0.001384233542160329"
$ws.Range("O7").Value = [double]"0.001384233542160329"
$ws.Range("P7").Value = [double]"-0.4654211338651546"
$ws.Range("Q7").Value = "[-1.0943686120613094, 0.16352634433100022]"
$ws.Range("R7").Value = [double]"0.1430861536172205"
$ws.Range("S7").Value = [double]"0.1430861536172205"
$ws.Range("T7").Value = [double]"64.56890744713847"
$ws.Range("U7").Value = "[47.75007084521431, 81.38774404906263]"
$ws.Range("V7").Value = [double]"8.477660795591646e-10"
$ws.Range("W7").Value = [double]"8.477660795591646e-10"
$ws.Range("X7").Value = [double]"1.900000000000041"
$ws.Range("Y7").Value = [double]"-0.6675675675675836"
$ws.Range("Z7").Value = [double]"4.467567567567666"

# Row 8
$ws.Range("B8").Value = [double]"0"
$ws.Range("F8").Value = [double]"25.8300000000006"
$ws.Range("H8").Value = [double]"0.00147592980018918"
$ws.Range("I8").Value = [double]"0.00147592980018918"
$ws.Range("L8").Value = [double]"46.0041110321439"
$ws.Range("M8").Value = "[14.282453690863036, 77.72576837342477]"
$ws.Range("N8").Value = [double]"0.005438269277630026"
$ws.Range("O8").Value = [double]"0.005438269277630026"
$ws.Range("P8").Value = [double]"-0.5283158816847697"
$ws.Range("Q8").Value = "[-1.2453160068283866, 0.18868424345884716]"
$ws.Range("R8").Value = [double]"0.1447599590504673"
$ws.Range("S8").Value = [double]"0.1447599590504673"
$ws.Range("T8").Value = [double]"70.60781175174131"
$ws.Range("U8").Value = "[53.671342438437534, 87.54428106504508]"
$ws.Range("V8").Value = [double]"9.18567444330165e-11"
$ws.Range("W8").Value = [double]"9.18567444330165e-11"
$ws.Range("X8").Value = [double]"2.171891891891942"
$ws.Range("Y8").Value = [double]"-0.775675675675696"
$ws.Range("Z8").Value = [double]"5.119459459459581"

# Row 9
$ws.Range("F9").Value = [double]"25.8300000000006"
$ws.Range("H9").Value = [double]"0.002620446759960982"
$ws.Range("I9").Value = [double]"0.002620446759960982"
$ws.Range("L9").Value = [double]"46.70352082680847"
$ws.Range("M9").Value = "[11.21843317652629, 82.18860847709065]"
$ws.Range("N9").Value = [double]"0.01104375419169523"
$ws.Range("O9").Value = [double]"0.01104375419169523"
$ws.Range("P9").Value = [double]"-1.044052813805617"
$ws.Range("Q9").Value = "[-1.8113687372049254, -0.2767368904063092]"
$ws.Range("R9").Value = [double]"0.008765952707644686"
$ws.Range("S9").Value = [double]"0.008765952707644686"
$ws.Range("T9").Value = [double]"64.75824925206408"
$ws.Range("U9").Value = "[46.06209111567657, 83.45440738845159]"
$ws.Range("V9").Value = [double]"1.104924418626752e-08"
$ws.Range("W9").Value = [double]"1.104924418626752e-08"
$ws.Range("X9").Value = [double]"4.292072072072173"
$ws.Range("Y9").Value = [double]"1.137657657657689"
$ws.Range("Z9").Value = [double]"7.446486486486657"

# Row 10
$ws.Range("F10").Value = [double]"25.8300000000006"
$ws.Range("H10").Value = [double]"3.746845155250256e-07"
$ws.Range("I10").Value = [double]"3.746845155250256e-07"
$ws.Range("L10").Value = [double]"58.01419727288999"
$ws.Range("M10").Value = "[36.26894501970554, 79.75944952607445]"
$ws.Range("N10").Value = [double]"2.628820495065298e-06"
$ws.Range("O10").Value = [double]"2.628820495065298e-06"
$ws.Range("P10").Value = [double]"-1.42142130072331"
$ws.Range("Q10").Value = "[-1.8491055858966954, -0.9937370155499243]"
$ws.Range("R10").Value = [double]"2.90178714390521e-08"
$ws.Range("S10").Value = [double]"2.90178714390521e-08"
$ws.Range("T10").Value = [double]"53.06029745944171"
$ws.Range("U10").Value = "[39.765950347090794, 66.35464457179263]"
$ws.Range("V10").Value = [double]"3.028282069550414e-10"
$ws.Range("W10").Value = [double]"3.028282069550414e-10"
$ws.Range("X10").Value = [double]"5.843423423423559"
$ws.Range("Y10").Value = [double]"4.085225225225319"
$ws.Range("Z10").Value = [double]"7.601621621621799"

# Row 11
$ws.Range("F11").Value = [double]"25.8300000000006"
$ws.Range("H11").Value = [double]"0.03384494813651173"
$ws.Range("I11").Value = [double]"0.03384494813651173"
$ws.Range("L11").Value = [double]"23.30666873605635"
$ws.Range("M11").Value = "[1.673048858897758, 44.94028861321494]"
$ws.Range("N11").Value = [double]"0.03533578663487158"
$ws.Range("O11").Value = [double]"0.03533578663487158"
$ws.Range("P11").Value = [double]"-1.773631888513156"
$ws.Range("Q11").Value = "[-3.107000542289004, -0.44026323473730855]"
$ws.Range("R11").Value = [double]"0.01027221178364379"
$ws.Range("S11").Value = [double]"0.01027221178364379"
$ws.Range("T11").Value = [double]"51.38300647005065"
$ws.Range("U11").Value = "[38.80202004840677, 63.963992891694524]"
$ws.Range("V11").Value = [double]"1.6202705843682e-10"
$ws.Range("W11").Value = [double]"1.6202705843682e-10"
$ws.Range("X11").Value = [double]"7.291351351351519"
$ws.Range("Y11").Value = [double]"1.809909909909953"
$ws.Range("Z11").Value = [double]"12.77279279279309"

# Row 12
$ws.Range("F12").Value = [double]"25.8300000000006"
$ws.Range("H12").Value = [double]"0.0004755474722017983"
$ws.Range("I12").Value = [double]"0.0004755474722017983"
$ws.Range("L12").Value = [double]"47.72923137970774"
$ws.Range("M12").Value = "[16.898474416660335, 78.55998834275515]"
$ws.Range("N12").Value = [double]"0.00317052415428809"
$ws.Range("O12").Value = [double]"0.00317052415428809"
$ws.Range("P12").Value = [double]"-2.352263568453619"
$ws.Range("Q12").Value = "[-2.9308952483940818, -1.7736318885131563]"
$ws.Range("R12").Value = [double]"1.840192442870148e-10"
$ws.Range("S12").Value = [double]"1.840192442870148e-10"
$ws.Range("T12").Value = [double]"59.89617793180389"
$ws.Range("U12").Value = "[44.28405853481361, 75.50829732879417]"
$ws.Range("V12").Value = [double]"8.626646064158194e-10"
$ws.Range("W12").Value = [double]"8.626646064158194e-10"
$ws.Range("X12").Value = [double]"9.670090090090316"
$ws.Range("Y12").Value = [double]"7.291351351351521"
$ws.Range("Z12").Value = [double]"12.04882882882911"

# Row 13
$ws.Range("F13").Value = [double]"25.8300000000006"
$ws.Range("H13").Value = [double]"0.001533953094238982"
$ws.Range("I13").Value = [double]"0.001533953094238982"
$ws.Range("L13").Value = [double]"44.22668731936476"
$ws.Range("M13").Value = "[15.48771016550782, 72.9656644732217]"
$ws.Range("N13").Value = [double]"0.003337961168694115"
$ws.Range("O13").Value = [double]"0.003337961168694115"
$ws.Range("P13").Value = [double]"-2.729632055371312"
$ws.Range("Q13").Value = "[-3.4466321805149276, -2.0126319302276956]"
$ws.Range("R13").Value = [double]"1.054208720319139e-09"
$ws.Range("S13").Value = [double]"1.054208720319139e-09"
$ws.Range("T13").Value = [double]"69.78581453209748"
$ws.Range("U13").Value = "[54.134422828051555, 85.4372062361434]"
$ws.Range("V13").Value = [double]"1.35047528715404e-11"
$ws.Range("W13").Value = [double]"1.35047528715404e-11"
$ws.Range("X13").Value = [double]"11.2214414414417"
$ws.Range("Y13").Value = [double]"8.273873873874066"
$ws.Range("Z13").Value = [double]"14.16900900900934"

# Row 14
$ws.Range("F14").Value = [double]"25.8300000000006"
$ws.Range("H14").Value = [double]"0.001111014688560452"
$ws.Range("I14").Value = [double]"0.001111014688560452"
$ws.Range("L14").Value = [double]"43.10029372547545"
$ws.Range("M14").Value = "[18.12099398928376, 68.07959346166713]"
$ws.Range("N14").Value = [double]"0.001141638123621336"
$ws.Range("O14").Value = [double]"0.001141638123621336"
$ws.Range("P14").Value = [double]"3.037816319687427"
$ws.Range("Q14").Value = "[2.3208161945438106, 3.7548164448310444]"
$ws.Range("R14").Value = [double]"5.844591477455197e-11"
$ws.Range("S14").Value = [double]"5.844591477455197e-11"
$ws.Range("T14").Value = [double]"59.85572443960125"
$ws.Range("U14").Value = "[44.85198693981965, 74.85946193938285]"
$ws.Range("V14").Value = [double]"3.065359077680796e-10"
$ws.Range("W14").Value = [double]"3.065359077680796e-10"
$ws.Range("X14").Value = [double]"13.34162162162193"
$ws.Range("Y14").Value = [double]"10.39405405405429"
$ws.Range("Z14").Value = [double]"16.28918918918957"
